$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9557.962
$ws.Range("I74").Value = 10021.125
$ws.Range("K74").Value = 10021.125
$ws.Range("M74").Value = -9085.125

$ws.Range("H77").Value = 9557.962
$ws.Range("I77").Value = 10021.125
$ws.Range("K77").Value = 50105.625
$ws.Range("M77").Value = -45425.625

$ws.Range("H103").Value = 91491.37
$ws.Range("I103").Value = 321.6
$ws.Range("K103").Value = 964.8000000000001
$ws.Range("M103").Value = -378.8000000000001

$ws.Range("H107").Value = 2080.2942
$ws.Range("I107").Value = 2422.3845
$ws.Range("K107").Value = 2422.3845
$ws.Range("M107").Value = -502.3845000000001

$ws.Range("H112").Value = 73065.28999999999
$ws.Range("J112").Value = 85094.25
$ws.Range("L112").Value = 255282.75
$ws.Range("N112").Value = -257498.75

$ws.Range("H113").Value = 5618.8276
$ws.Range("I113").Value = 4140.25
$ws.Range("K113").Value = 4140.25
$ws.Range("M113").Value = -886.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 405981.72
$ws.Range("I32").Value = 405981.72
$ws.Range("K32").Value = 405981.72
$ws.Range("M32").Value = -405694.72

$ws.Range("H45").Value = 172497
$ws.Range("I45").Value = 172497
$ws.Range("K45").Value = 172497
$ws.Range("M45").Value = -172120

$ws.Range("H61").Value = 4214202
$ws.Range("I61").Value = 16478.111
$ws.Range("K61").Value = 16478.111
$ws.Range("M61").Value = -16266.111

$ws.Range("H64").Value = 30080
$ws.Range("I64").Value = 30080
$ws.Range("J64").Value = 30080
$ws.Range("K64").Value = 30080
$ws.Range("L64").Value = 30080
$ws.Range("M64").Value = -29832
$ws.Range("N64").Value = -30576

$ws.Range("H67").Value = 30080
$ws.Range("I67").Value = 30080
$ws.Range("J67").Value = 30080
$ws.Range("K67").Value = 30080
$ws.Range("L67").Value = 30080
$ws.Range("M67").Value = -29222
$ws.Range("N67").Value = -31796

$ws.Range("H68").Value = 37567.25
$ws.Range("J68").Value = 37567.25
$ws.Range("L68").Value = 37567.25
$ws.Range("N68").Value = -39189.25

$ws.Range("H71").Value = 37567.25
$ws.Range("J71").Value = 37567.25
$ws.Range("L71").Value = 112701.75
$ws.Range("N71").Value = -120813.75

$ws.Range("H75").Value = 82533
$ws.Range("J75").Value = 82533
$ws.Range("L75").Value = 82533
$ws.Range("N75").Value = -84281

$ws.Range("H78").Value = 82533
$ws.Range("J78").Value = 82533
$ws.Range("L78").Value = 247599
$ws.Range("N78").Value = -256335

$ws.Range("H97").Value = 6222.263
$ws.Range("I97").Value = 13622.375
$ws.Range("J97").Value = 840.36365
$ws.Range("K97").Value = 13622.375
$ws.Range("L97").Value = 840.36365
$ws.Range("M97").Value = -13126.375
$ws.Range("N97").Value = -1832.36365

$ws.Range("H102").Value = 2810.1875
$ws.Range("I102").Value = 2078.8333
$ws.Range("J102").Value = 5004.25
$ws.Range("K102").Value = 2078.8333
$ws.Range("L102").Value = 5004.25
$ws.Range("M102").Value = -456.8332999999998
$ws.Range("N102").Value = -8248.25

$ws.Range("H136").Value = 4214202
$ws.Range("I136").Value = 16478.111
$ws.Range("K136").Value = 49434.333
$ws.Range("M136").Value = -46884.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1763.5
$ws.Range("I94").Value = 1858.7778
$ws.Range("J94").Value = 1592
$ws.Range("K94").Value = 1858.7778
$ws.Range("L94").Value = 1592
$ws.Range("M94").Value = -1407.7778
$ws.Range("N94").Value = -2494

$ws.Range("H96").Value = 16777.25
$ws.Range("I96").Value = 16777.25
$ws.Range("K96").Value = 16777.25
$ws.Range("M96").Value = -14031.25

$ws.Range("H99").Value = 2837.75
$ws.Range("I99").Value = 3321.8
$ws.Range("J99").Value = 2710.3684
$ws.Range("K99").Value = 3321.8
$ws.Range("L99").Value = 2710.3684
$ws.Range("M99").Value = -1823.8
$ws.Range("N99").Value = -5706.368399999999

$ws.Range("H134").Value = 29034532
$ws.Range("I134").Value = 2088.4644
$ws.Range("K134").Value = 6265.3932
$ws.Range("M134").Value = -3730.3932

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3349.6365
$ws.Range("J31").Value = 4694.227
$ws.Range("L31").Value = 4694.227
$ws.Range("N31").Value = -5284.227

$ws.Range("H34").Value = 3349.6365
$ws.Range("J34").Value = 4694.227
$ws.Range("L34").Value = 4694.227
$ws.Range("N34").Value = -5098.227

$ws.Range("H86").Value = 164404
$ws.Range("I86").Value = 504557.5
$ws.Range("J86").Value = 28342.6
$ws.Range("K86").Value = 504557.5
$ws.Range("L86").Value = 28342.6
$ws.Range("M86").Value = -503434.5
$ws.Range("N86").Value = -30588.6

$ws.Range("H89").Value = 164404
$ws.Range("I89").Value = 504557.5
$ws.Range("J89").Value = 28342.6
$ws.Range("K89").Value = 2522787.5
$ws.Range("L89").Value = 141713
$ws.Range("M89").Value = -2517171.5
$ws.Range("N89").Value = -152945

$ws.Range("H94").Value = 1242
$ws.Range("I94").Value = 1097
$ws.Range("J94").Value = 1416
$ws.Range("K94").Value = 1097
$ws.Range("L94").Value = 1416
$ws.Range("M94").Value = -646
$ws.Range("N94").Value = -2318

$ws.Range("H133").Value = 117496
$ws.Range("J133").Value = 117496
$ws.Range("L133").Value = 117496
$ws.Range("N133").Value = -122556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 125037.625
$ws.Range("J23").Value = 166703.5
$ws.Range("L23").Value = 500110.5
$ws.Range("N23").Value = -500580.5

$ws.Range("H121").Value = 2266
$ws.Range("I121").Value = 1350
$ws.Range("J121").Value = 2998.8
$ws.Range("K121").Value = 4050
$ws.Range("L121").Value = 8996.400000000001
$ws.Range("M121").Value = -2740
$ws.Range("N121").Value = -11616.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 21497.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 21497.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 21497.5
$ws.Range("M54").Value = $null
$ws.Range("N54").Value = -22277.5

$ws.Range("H107").Value = 71714.57000000001
$ws.Range("I107").Value = 143224.42
$ws.Range("J107").Value = 204.71428
$ws.Range("K107").Value = 143224.42
$ws.Range("L107").Value = 204.71428
$ws.Range("M107").Value = -141304.42
$ws.Range("N107").Value = -4044.71428

$ws.Range("H126").Value = 2068.875
$ws.Range("I126").Value = 1961.7693
$ws.Range("K126").Value = 5885.3079
$ws.Range("M126").Value = -3415.3079

$ws.Range("H135").Value = 107852
$ws.Range("J135").Value = 107852
$ws.Range("L135").Value = 107852
$ws.Range("N135").Value = -117992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2721.4285
$ws.Range("I7").Value = 2369.7778
$ws.Range("K7").Value = 2369.7778
$ws.Range("M7").Value = -2257.7778

$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3340

$ws.Range("H22").Value = 1810
$ws.Range("J22").Value = 1992.5
$ws.Range("L22").Value = 1992.5
$ws.Range("N22").Value = -2582.5

$ws.Range("H27").Value = 1810
$ws.Range("J27").Value = 1992.5
$ws.Range("L27").Value = 1992.5
$ws.Range("N27").Value = -2206.5

$ws.Range("H100").Value = 1651.0834
$ws.Range("I100").Value = 1606.5714
$ws.Range("J100").Value = 1713.4
$ws.Range("K100").Value = 1606.5714
$ws.Range("L100").Value = 1713.4
$ws.Range("M100").Value = -1065.5714
$ws.Range("N100").Value = -2795.4

$ws.Range("H122").Value = 4858.356
$ws.Range("I122").Value = 3972.7942
$ws.Range("K122").Value = 11918.3826
$ws.Range("M122").Value = -9468.382599999999

$ws.Range("H126").Value = 2721.4285
$ws.Range("I126").Value = 2369.7778
$ws.Range("K126").Value = 7109.3334
$ws.Range("M126").Value = -4639.3334

$ws.Range("H136").Value = 3599.5417
$ws.Range("I136").Value = 1679
$ws.Range("J136").Value = 4751.8667
$ws.Range("K136").Value = 5037
$ws.Range("L136").Value = 14255.6001
$ws.Range("M136").Value = -2487
$ws.Range("N136").Value = -19355.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8420.817999999999
$ws.Range("I62").Value = 8156
$ws.Range("J62").Value = 8520.125
$ws.Range("K62").Value = 8156
$ws.Range("L62").Value = 8520.125
$ws.Range("M62").Value = -7532
$ws.Range("N62").Value = -9768.125

$ws.Range("H65").Value = 8420.817999999999
$ws.Range("I65").Value = 8156
$ws.Range("J65").Value = 8520.125
$ws.Range("K65").Value = 40780
$ws.Range("L65").Value = 42600.625
$ws.Range("M65").Value = -37660
$ws.Range("N65").Value = -48840.625

$ws.Range("H81").Value = 3791.4167
$ws.Range("I81").Value = 2583.111
$ws.Range("J81").Value = 7416.3335
$ws.Range("K81").Value = 5166.222
$ws.Range("L81").Value = 14832.667
$ws.Range("M81").Value = -4105.222
$ws.Range("N81").Value = -16954.667

$ws.Range("H84").Value = 3791.4167
$ws.Range("I84").Value = 2583.111
$ws.Range("J84").Value = 7416.3335
$ws.Range("K84").Value = 25831.11
$ws.Range("L84").Value = 74163.33499999999
$ws.Range("M84").Value = -20527.11
$ws.Range("N84").Value = -84771.33499999999

$ws.Range("H100").Value = 732.82355
$ws.Range("I100").Value = 673.7692
$ws.Range("J100").Value = 924.75
$ws.Range("K100").Value = 1347.5384
$ws.Range("L100").Value = 1849.5
$ws.Range("M100").Value = -806.5383999999999
$ws.Range("N100").Value = -2931.5

$ws.Range("H107").Value = 40000388
$ws.Range("I107").Value = 362.58823
$ws.Range("K107").Value = 1087.76469
$ws.Range("M107").Value = 832.23531

$ws.Range("H113").Value = 353.375
$ws.Range("I113").Value = 343.16666
$ws.Range("J113").Value = 363.58334
$ws.Range("K113").Value = 1029.49998
$ws.Range("L113").Value = 1090.75002
$ws.Range("M113").Value = 1140.50002
$ws.Range("N113").Value = -5430.750019999999
